$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.030222437722262
$ws.Range("D2").Value = 1.035566803709351
$ws.Range("E2").Value = 1.029952765763676
$ws.Range("F2").Value = 1.041996612835217
$ws.Range("I2").Value = 1.037866156299475
$ws.Range("J2").Value = 1.035364874375759
$ws.Range("K2").Value = 1.038363039763212
$ws.Range("L2").Value = 1.03276518387113
$ws.Range("M2").Value = 1.044774547989205
$ws.Range("N2").Value = 1.036835211255401
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.031104783179362
$ws.Range("D3").Value = 1.036265810512665
$ws.Range("E3").Value = 1.030700757619069
$ws.Range("F3").Value = 1.044226453738626
$ws.Range("I3").Value = 1.038177146327951
$ws.Range("J3").Value = 1.035888925670716
$ws.Range("K3").Value = 1.038871761577307
$ws.Range("L3").Value = 1.033321591033351
$ws.Range("M3").Value = 1.046811413919474
$ws.Range("N3").Value = 1.037360006763289
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.031675734904748
$ws.Range("D4").Value = 1.036717932095541
$ws.Range("E4").Value = 1.031185141763024
$ws.Range("F4").Value = 1.045663758802845
$ws.Range("I4").Value = 1.038376659397961
$ws.Range("J4").Value = 1.036227404932204
$ws.Range("K4").Value = 1.039200082571542
$ws.Range("L4").Value = 1.033681344037096
$ws.Range("M4").Value = 1.048123528012336
$ws.Range("N4").Value = 1.037698966704152
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.031915767070637
$ws.Range("D5").Value = 1.036907960724953
$ws.Range("E5").Value = 1.031388868545573
$ws.Range("F5").Value = 1.046266705444914
$ws.Range("I5").Value = 1.038460124979242
$ws.Range("J5").Value = 1.036369554426658
$ws.Range("K5").Value = 1.039337904545286
$ws.Range("L5").Value = 1.033832517495369
$ws.Range("M5").Value = 1.048673764758163
$ws.Range("N5").Value = 1.037841318067193
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.031956069777289
$ws.Range("D6").Value = 1.036939864841439
$ws.Range("E6").Value = 1.031423080542154
$ws.Range("F6").Value = 1.046367867754424
$ws.Range("I6").Value = 1.038474115238981
$ws.Range("J6").Value = 1.036393413335688
$ws.Range("K6").Value = 1.03936103351004
$ws.Range("L6").Value = 1.033857896279196
$ws.Range("M6").Value = 1.048766072111573
$ws.Range("N6").Value = 1.03786521085861
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.031678942209868
$ws.Range("D7").Value = 1.036720471437008
$ws.Range("E7").Value = 1.031187863610167
$ws.Range("F7").Value = 1.04567182045982
$ws.Range("I7").Value = 1.038377776277055
$ws.Range("J7").Value = 1.036229304918459
$ws.Range("K7").Value = 1.039201924956866
$ws.Range("L7").Value = 1.033683364286097
$ws.Range("M7").Value = 1.048130885672652
$ws.Range("N7").Value = 1.037700869388606
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.030520627610024
$ws.Range("D8").Value = 1.035803074394884
$ws.Range("E8").Value = 1.030205473677595
$ws.Range("F8").Value = 1.042751367114092
$ws.Range("I8").Value = 1.037971613615539
$ws.Range("J8").Value = 1.035542108030322
$ws.Range("K8").Value = 1.03853514239921
$ws.Range("L8").Value = 1.032953282329261
$ws.Range("M8").Value = 1.045464151059948
$ws.Range("N8").Value = 1.037012696602065
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.028479625283587
$ws.Range("D9").Value = 1.034185084458977
$ws.Range("E9").Value = 1.028477316000476
$ws.Range("F9").Value = 1.037561166417176
$ws.Range("I9").Value = 1.037242667807848
$ws.Range("J9").Value = 1.034326429353196
$ws.Range("K9").Value = 1.037353594113317
$ws.Range("L9").Value = 1.031664630618117
$ws.Range("M9").Value = 1.040718659865583
$ws.Range("N9").Value = 1.035795291521778
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.027118988314643
$ws.Range("D10").Value = 1.033105436423315
$ws.Range("E10").Value = 1.02732718592052
$ws.Range("F10").Value = 1.034069378622976
$ws.Range("I10").Value = 1.036747696795409
$ws.Range("J10").Value = 1.033512745164819
$ws.Range("K10").Value = 1.036561408910787
$ws.Range("L10").Value = 1.030804060341746
$ws.Range("M10").Value = 1.037521889291218
$ws.Range("N10").Value = 1.034980451808522
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.026529816079208
$ws.Range("D11").Value = 1.032637693576621
$ws.Range("E11").Value = 1.026829633784706
$ws.Range("F11").Value = 1.032549420474693
$ws.Range("I11").Value = 1.03653120743088
$ws.Range("J11").Value = 1.033159634821531
$ws.Range("K11").Value = 1.036217306680817
$ws.Range("L11").Value = 1.03043106969008
$ws.Range("M11").Value = 1.036129363216092
$ws.Range("N11").Value = 1.034626840008044
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.026310969265036
$ws.Range("D12").Value = 1.032463914877766
$ws.Range("E12").Value = 1.026644889888234
$ws.Range("F12").Value = 1.031983599191064
$ws.Range("I12").Value = 1.036450466321165
$ws.Range("J12").Value = 1.033028355842255
$ws.Range("K12").Value = 1.036089328164532
$ws.Range("L12").Value = 1.030292469961573
$ws.Range("M12").Value = 1.0356108316845
$ws.Range("N12").Value = 1.03449537459756
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.026357912755411
$ws.Range("D13").Value = 1.032501192741234
$ws.Range("E13").Value = 1.026684514957815
$ws.Range("F13").Value = 1.032105026600978
$ws.Range("I13").Value = 1.03646780041199
$ws.Range("J13").Value = 1.033056520987449
$ws.Range("K13").Value = 1.036116787415648
$ws.Range("L13").Value = 1.030322202545227
$ws.Range("M13").Value = 1.035722117237624
$ws.Range("N13").Value = 1.03452357974049
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.026511726179651
$ws.Range("D14").Value = 1.032623329759395
$ws.Range("E14").Value = 1.026814361388175
$ws.Range("F14").Value = 1.032502675067381
$ws.Range("I14").Value = 1.036524540036042
$ws.Range("J14").Value = 1.033148785674376
$ws.Range("K14").Value = 1.036206731282049
$ws.Range("L14").Value = 1.030419614106485
$ws.Range("M14").Value = 1.036086527680019
$ws.Range("N14").Value = 1.034615975453855
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.026606495412675
$ws.Range("D15").Value = 1.032698577320407
$ws.Range("E15").Value = 1.026894373210202
$ws.Range("F15").Value = 1.032747513720395
$ws.Range("I15").Value = 1.036559455749532
$ws.Range("J15").Value = 1.033205617318466
$ws.Range("K15").Value = 1.03626212694645
$ws.Range("L15").Value = 1.030479625358197
$ws.Range("M15").Value = 1.036310881515793
$ws.Range("N15").Value = 1.0346728878054
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.027158089503126
$ws.Range("D16").Value = 1.033136473672001
$ws.Range("E16").Value = 1.02736021652517
$ws.Range("F16").Value = 1.034170081429378
$ws.Range("I16").Value = 1.036762018704466
$ws.Range("J16").Value = 1.033536163428879
$ws.Range("K16").Value = 1.036584222965589
$ws.Range("L16").Value = 1.03082880690877
$ws.Range("M16").Value = 1.037614128502189
$ws.Range("N16").Value = 1.035003903329203
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.02750408728898
$ws.Range("D17").Value = 1.033411087432297
$ws.Range("E17").Value = 1.02765255119001
$ws.Range("F17").Value = 1.035060252715859
$ws.Range("I17").Value = 1.036888500246579
$ws.Range("J17").Value = 1.033743296871299
$ws.Range("K17").Value = 1.036785975057657
$ws.Range("L17").Value = 1.031047742997567
$ws.Range("M17").Value = 1.038429371345465
$ws.Range("N17").Value = 1.035211330924869
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.027705901288537
$ws.Range("D18").Value = 1.033571241116764
$ws.Range("E18").Value = 1.02782310962943
$ws.Range("F18").Value = 1.035578706328582
$ws.Range("I18").Value = 1.036962066184887
$ws.Range("J18").Value = 1.033864039137796
$ws.Range("K18").Value = 1.036903549444432
$ws.Range("L18").Value = 1.031175410211334
$ws.Range("M18").Value = 1.038904090143104
$ws.Range("N18").Value = 1.035332244659234
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.027774714561932
$ws.Range("D19").Value = 1.033625845360857
$ws.Range("E19").Value = 1.027881273217192
$ws.Range("F19").Value = 1.035755356524139
$ws.Range("I19").Value = 1.036987114952991
$ws.Range("J19").Value = 1.033905196408267
$ws.Range("K19").Value = 1.036943621635421
$ws.Range("L19").Value = 1.031218935597107
$ws.Range("M19").Value = 1.039065822671825
$ws.Range("N19").Value = 1.035373460377751
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.027466965075378
$ws.Range("D20").Value = 1.033381626450879
$ws.Range("E20").Value = 1.027621181838392
$ws.Range("F20").Value = 1.034964825515575
$ws.Range("I20").Value = 1.036874951567412
$ws.Range("J20").Value = 1.033721081194761
$ws.Range("K20").Value = 1.036764339753966
$ws.Range("L20").Value = 1.031024256800398
$ws.Range("M20").Value = 1.03834198637473
$ws.Range("N20").Value = 1.035189083699523
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.026466432022845
$ws.Range("D21").Value = 1.032587364513139
$ws.Range("E21").Value = 1.026776122944953
$ws.Range("F21").Value = 1.032385612096494
$ws.Range("I21").Value = 1.036507840684979
$ws.Range("J21").Value = 1.033121619290835
$ws.Range("K21").Value = 1.036180249582319
$ws.Range("L21").Value = 1.030390930336121
$ws.Range("M21").Value = 1.035979253676582
$ws.Range("N21").Value = 1.034588770490934
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.025837345064312
$ws.Range("D22").Value = 1.032087759134042
$ws.Range("E22").Value = 1.02624520063322
$ws.Range("F22").Value = 1.030756756710103
$ws.Range("I22").Value = 1.036275128571855
$ws.Range("J22").Value = 1.032744030162467
$ws.Range("K22").Value = 1.035812061564685
$ws.Range("L22").Value = 1.029992417827519
$ws.Range("M22").Value = 1.034486254005086
$ws.Range("N22").Value = 1.034210645142693
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.026170837420062
$ws.Range("D23").Value = 1.032352630694858
$ws.Range("E23").Value = 1.026526614798582
$ws.Range("F23").Value = 1.031620940130176
$ws.Range("I23").Value = 1.036398674050121
$ws.Range("J23").Value = 1.03294426242464
$ws.Range("K23").Value = 1.03600733522381
$ws.Range("L23").Value = 1.030203706918117
$ws.Range("M23").Value = 1.035278440922098
$ws.Range("N23").Value = 1.034411161757646
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.027483739008259
$ws.Range("D24").Value = 1.033394938674364
$ws.Range("E24").Value = 1.027635356160768
$ws.Range("F24").Value = 1.035007947332012
$ws.Range("I24").Value = 1.036881074276991
$ws.Range("J24").Value = 1.033731119735241
$ws.Range("K24").Value = 1.036774116138614
$ws.Range("L24").Value = 1.031034869308844
$ws.Range("M24").Value = 1.038381474347611
$ws.Range("N24").Value = 1.035199136495882
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.029007265240068
$ws.Range("D25").Value = 1.034603545160406
$ws.Range("E25").Value = 1.028923737886028
$ws.Range("F25").Value = 1.038908382505991
$ws.Range("I25").Value = 1.037432697234461
$ws.Range("J25").Value = 1.034641278135998
$ws.Range("K25").Value = 1.037659839314962
$ws.Range("L25").Value = 1.031998035281274
$ws.Range("M25").Value = 1.041951174963077
$ws.Range("N25").Value = 1.036110587425965

Write-Output "Applied 264 cell updates to vm_pu sheet (380 kV case)."
